$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 93-95 (MarketObjects column) ---
$ws.Range("B93").Value = "['USD.SOFR.CSA_USD', 'BTCUSD.SPOT', 'BTC.FUNDING.CSA_USD', 'BTCUSD.QPROBABILITY']"
$ws.Range("B94").Value = "['BTCUSD.SPOT']"
$ws.Range("B95").Value = "['BTCUSD.SPOT']"

# --- Append new rows 96-100, keeping the Date column as plain text ---
# Force text number format before assigning so Excel doesn't auto-detect
# the "YYYY-MM-DD" strings as dates, then clear the formatting override so
# the cells end up styleless (matching the rest of the column).
$ws.Range("A96:A100").NumberFormat = "@"

$ws.Range("A96").Value = "2025-09-15"
$ws.Range("B96").Value = "['USD.SOFR.CSA_USD', 'BTCUSD.SPOT', 'BTC.FUNDING.CSA_USD', 'BTCUSD.QPROBABILITY']"

$ws.Range("A97").Value = "2025-09-16"
$ws.Range("B97").Value = "['USD.SOFR.CSA_USD', 'BTCUSD.SPOT', 'BTC.FUNDING.CSA_USD', 'BTCUSD.QPROBABILITY']"

$ws.Range("A98").Value = "2025-09-17"
$ws.Range("B98").Value = "['USD.SOFR.CSA_USD', 'BTCUSD.SPOT', 'BTC.FUNDING.CSA_USD', 'BTCUSD.QPROBABILITY']"

$ws.Range("A99").Value = "2025-09-18"
$ws.Range("B99").Value = "['USD.SOFR.CSA_USD', 'BTCUSD.SPOT']"

$ws.Range("A100").Value = "2025-09-19"
$ws.Range("B100").Value = "['USD.SOFR.CSA_USD']"

$ws.Range("A96:A100").ClearFormats()
